# HR1 case study updated.
#
# - "DownActivation, 2020, Winter" selection moves to A7 (was B1:Y1).
# - "UpActivation, 2020, Winter" becomes the active/selected tab, with its
#   selection moved to A7 (was S11); this naturally clears tabSelected on
#   the previously-active "Investment Cost" sheet.
# - "UpActivation, 2020, Winter" column A (rows 2-6) is reindexed to start
#   at 0 instead of 1.

$wb = $excel.ActiveWorkbook

$wsDown = $wb.Worksheets.Item("DownActivation, 2020, Winter")
$wsUp   = $wb.Worksheets.Item("UpActivation, 2020, Winter")

# Update the data values in "UpActivation, 2020, Winter" column A.
$wsUp.Range("A2").Value = 0
$wsUp.Range("A3").Value = 1
$wsUp.Range("A4").Value = 2
$wsUp.Range("A5").Value = 3
$wsUp.Range("A6").Value = 4

# Move the selection on "DownActivation, 2020, Winter" to A7 (not the
# active tab, so select without leaving it active last).
$wsDown.Activate()
$wsDown.Range("A7").Select()

# Finally, activate "UpActivation, 2020, Winter" and move its selection to
# A7 as well; being the last sheet activated, it becomes the workbook's
# active tab (tabSelected="1" / activeTab="4").
$wsUp.Activate()
$wsUp.Range("A7").Select()
